$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'youth black knee pads'
$ws.Cells.Item(2, 1).Value = 'basketball clothes for boys'
$ws.Cells.Item(3, 1).Value = 'men capri'
$ws.Cells.Item(4, 1).Value = 'black baseball pants adult small'
$ws.Cells.Item(5, 1).Value = 'under pants for men'
$ws.Cells.Item(6, 1).Value = 'by design knee pads'
$ws.Cells.Item(7, 1).Value = 'baseball knee high pants mens'
$ws.Cells.Item(8, 1).Value = 'boys youth compression tights'
$ws.Cells.Item(9, 1).Value = 'mens small running tights'
$ws.Cells.Item(10, 1).Value = 'youth large softball pants'
$ws.Cells.Item(11, 1).Value = 'small basketballs'
$ws.Cells.Item(12, 1).Value = 'cycling capri'
$ws.Cells.Item(13, 1).Value = 'mens black basketball shorts'
$ws.Cells.Item(14, 1).Value = 'knee protector for running'
$ws.Cells.Item(15, 1).Value = 'knees protection'
$ws.Cells.Item(16, 1).Value = 'dry fit baseball pants'
$ws.Cells.Item(17, 1).Value = 'girls black baseball pants'
$ws.Cells.Item(18, 1).Value = 'mens baseball sliding shorts'
$ws.Cells.Item(19, 1).Value = 'wrestling equipment'
$ws.Cells.Item(20, 1).Value = 'the knee pads'
$ws.Cells.Item(21, 1).Value = 'football compression'
$ws.Cells.Item(22, 1).Value = 'soccer guards for men'
$ws.Cells.Item(23, 1).Value = 'spandex leggings for boys'
$ws.Cells.Item(24, 1).Value = 'basketball pants men'
$ws.Cells.Item(25, 1).Value = 'hunting knee pads'
$ws.Cells.Item(26, 1).Value = 'combat pants knee pads'
$ws.Cells.Item(27, 1).Value = 'nike compression basketball pants'
$ws.Cells.Item(28, 1).Value = 'seamless capri leggings'
$ws.Cells.Item(29, 1).Value = 'basketball knee pads adidas'
$ws.Cells.Item(30, 1).Value = 'women black leggings'
$ws.Cells.Item(31, 1).Value = 'red basketball knee pads'
$ws.Cells.Item(32, 1).Value = 'underarmour mens leggings'
$ws.Cells.Item(33, 1).Value = 'red nike compression pants men'
$ws.Cells.Item(34, 1).Value = 'men compression pants long'
$ws.Cells.Item(35, 1).Value = 'knee tights'
$ws.Cells.Item(36, 1).Value = 'knee protector basketball'
$ws.Cells.Item(37, 1).Value = 'youth football pants with pads'
$ws.Cells.Item(38, 1).Value = 'little boys compression tights'
$ws.Cells.Item(39, 1).Value = 'wrestling tights for boys'
$ws.Cells.Item(40, 1).Value = 'youth boys baseball pants'
$ws.Cells.Item(41, 1).Value = 'padded basketball tights'
$ws.Cells.Item(42, 1).Value = 'basketball clothes youth'
$ws.Cells.Item(43, 1).Value = 'leggings boys'
$ws.Cells.Item(44, 1).Value = 'youth compression pants'
$ws.Cells.Item(45, 1).Value = 'knee pad youth'
$ws.Cells.Item(46, 1).Value = 'best knee pads for volleyball'
$ws.Cells.Item(47, 1).Value = 'athletic pants for men big and tall'
$ws.Cells.Item(48, 1).Value = 'leg compression basketball'
$ws.Cells.Item(49, 1).Value = 'compression shorts baseball'
$ws.Cells.Item(50, 1).Value = 'pants for men sports'
$ws.Cells.Item(51, 1).Value = 'knee pads large'
$ws.Cells.Item(52, 1).Value = 'gym tights for men'
$ws.Cells.Item(53, 1).Value = 'soccer gear for men'
$ws.Cells.Item(54, 1).Value = 'large tall athletic pants men'
$ws.Cells.Item(55, 1).Value = 'youth football pads'
$ws.Cells.Item(56, 1).Value = 'padded compression shorts'
$ws.Cells.Item(57, 1).Value = 'youth girls softball pants black'
$ws.Cells.Item(58, 1).Value = 'athletic pants for men'
$ws.Cells.Item(59, 1).Value = 'patella knee'
$ws.Cells.Item(60, 1).Value = 'youth softball pants'
$ws.Cells.Item(61, 1).Value = 'basketball knee sleeves'
$ws.Cells.Item(62, 1).Value = 'knee pads for crossfit'
$ws.Cells.Item(63, 1).Value = 'mens nike dri fit compression pants'
$ws.Cells.Item(64, 1).Value = 'volleyball youth knee pads'
$ws.Cells.Item(65, 1).Value = 'troll knee pads'
$ws.Cells.Item(66, 1).Value = 'elastic knee pads'
$ws.Cells.Item(67, 1).Value = 'men leggings adidas'
$ws.Cells.Item(68, 1).Value = 'snowmobile knee pads'
$ws.Cells.Item(69, 1).Value = 'knee pads military'
$ws.Cells.Item(70, 1).Value = 'knee pads sleeve'
$ws.Cells.Item(71, 1).Value = 'knee pads mma'
$ws.Cells.Item(72, 1).Value = 'under armour compression pants men'
$ws.Cells.Item(73, 1).Value = 'mcdavid knee pads basketball'
$ws.Cells.Item(74, 1).Value = 'black leggings xsmall'
$ws.Cells.Item(75, 1).Value = 'tesla compression pants'
$ws.Cells.Item(76, 1).Value = 'fox knee pads'
$ws.Cells.Item(77, 1).Value = 'biking knee pads'
$ws.Cells.Item(78, 1).Value = 'bmx knee pads'
$ws.Cells.Item(79, 1).Value = 'mcdavid knee pad'
$ws.Cells.Item(80, 1).Value = 'pant with knee pads'
$ws.Cells.Item(81, 1).Value = 'green knee pads'
$ws.Cells.Item(82, 1).Value = 'knee pads tights'
$ws.Cells.Item(83, 1).Value = 'warm compression pants'
$ws.Cells.Item(84, 1).Value = 'dancer knee pads'
$ws.Cells.Item(85, 1).Value = 'navy compression pants'
$ws.Cells.Item(86, 1).Value = 'mizuno knee pad'
$ws.Cells.Item(87, 1).Value = 'mens pants with knee pads'
$ws.Cells.Item(88, 1).Value = 'black athletic capri'
$ws.Cells.Item(89, 1).Value = 'tactical knee pad'
$ws.Cells.Item(90, 1).Value = 'men compression pants blue'
$ws.Cells.Item(91, 1).Value = 'leggins for men sport'
$ws.Cells.Item(92, 1).Value = 'knee pads ski'
$ws.Cells.Item(93, 1).Value = 'knee pad snowboard'
$ws.Cells.Item(94, 1).Value = 'knee pads mcdavid basketball'
$ws.Cells.Item(95, 1).Value = 'pantalones con rodilleras'
$ws.Cells.Item(96, 1).Value = 'pantalon con rodilleras'
$ws.Cells.Item(97, 1).Value = 'cat knee pad pants'
$ws.Cells.Item(98, 1).Value = 'knee padded pants men'
$ws.Cells.Item(99, 1).Value = 'knee pad pants men'
$ws.Cells.Item(100, 1).Value = 'caterpillar knee pad pants'
